# Add an "Entity" column to the phishing-samples dataset.
# New column is inserted before the existing "url" column (old column J),
# pushing url to column K, and the new J column is populated with the
# entity/brand that each phishing sample impersonated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J; this shifts the existing J (url) column to K
# and carries its formatting (bestFit width, etc.) along with it.
$ws.Columns("J:J").Insert()

# Header for the new column.
$ws.Range("J1").Value = "Entity"

# New column shouldn't inherit the old bestFit/url width - give it its own.
$ws.Columns("J:J").ColumnWidth = 15.877604166666668

# Fix a small typo in the description of row 19 while we're at it
# ("2facode" -> "2fa code").
$ws.Range("I19").Value = "request for 2fa code"

# Populate the Entity column for every data row.
$entities = @{
    "2"  = "DHL"
    "3"  = "MaltaPost"
    "4"  = "HSBC"
    "5"  = "MaltaPost"
    "6"  = "BOV,HSBC"
    "7"  = "Government"
    "8"  = "MaltaPost"
    "9"  = "MaltaPost"
    "10" = "MaltaPost"
    "11" = "MaltaPost"
    "12" = "MaltaPost"
    "13" = "MaltaPost"
    "14" = "MaltaPost"
    "15" = "MaltaPost"
    "16" = "MaltaPost"
    "17" = "None"
    "18" = "BOV"
    "19" = "None"
    "20" = "CentralBank"
    "21" = "Melita"
    "22" = "MTA"
    "23" = "MCAST"
    "24" = "ChinaUniversity"
    "25" = "LIDL"
    "26" = "Toyota"
    "27" = "Netflix"
    "28" = "MeDirect"
    "29" = "DHL"
    "30" = "None"
    "31" = "None"
    "32" = "BOV"
}

foreach ($row in $entities.Keys) {
    $ws.Range("J$row").Value = $entities[$row]
}

# Leave the selection where the author ended up after entering the data.
$ws.Range("J33").Select()
